$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.198019801980198
$ws.Range("C2").Value = 0.528052805280528
$ws.Range("J2").Value = 0.0264026402640264
$ws.Range("O2").Value = 0.0033003300330033
$ws.Range("P2").Value = 0.1386138613861386
$ws.Range("S2").Value = 0.1056105610561056
$ws.Range("B3").Value = 0.02366863905325444
$ws.Range("C3").Value = 0.04142011834319527
$ws.Range("J3").Value = 0.02366863905325444
$ws.Range("P3").Value = 0.7218934911242604
$ws.Range("S3").Value = 0.1893491124260355
$ws.Range("J4").Value = 0.02127659574468085
$ws.Range("P4").Value = 0.8297872340425532
$ws.Range("S4").Value = 0.148936170212766
$ws.Range("B6").Value = 0.07511737089201878
$ws.Range("D6").Value = 0.004694835680751174
$ws.Range("F6").Value = 0.05633802816901409
$ws.Range("J6").Value = 0.2676056338028169
$ws.Range("O6").Value = 0.02347417840375587
$ws.Range("Q6").Value = 0.1643192488262911
$ws.Range("R6").Value = 0.07042253521126761
$ws.Range("S6").Value = 0.3380281690140845
$ws.Range("B7").Value = 0.1129943502824859
$ws.Range("D7").Value = 0.005649717514124294
$ws.Range("F7").Value = 0.06779661016949153
$ws.Range("J7").Value = 0.1638418079096045
$ws.Range("O7").Value = 0.02824858757062147
$ws.Range("Q7").Value = 0.192090395480226
$ws.Range("R7").Value = 0.05649717514124294
$ws.Range("S7").Value = 0.3728813559322034
$ws.Range("B8").Value = 0.09192825112107623
$ws.Range("D8").Value = 0.02242152466367713
$ws.Range("F8").Value = 0.05381165919282511
$ws.Range("J8").Value = 0.1255605381165919
$ws.Range("O8").Value = 0.02690582959641256
$ws.Range("Q8").Value = 0.1704035874439462
$ws.Range("R8").Value = 0.1121076233183857
$ws.Range("S8").Value = 0.3968609865470852
$ws.Range("B9").Value = 0.1165644171779141
$ws.Range("D9").Value = 0.01840490797546012
$ws.Range("F9").Value = 0.049079754601227
$ws.Range("J9").Value = 0.2024539877300613
$ws.Range("O9").Value = 0.01226993865030675
$ws.Range("Q9").Value = 0.1779141104294479
$ws.Range("R9").Value = 0.0736196319018405
$ws.Range("S9").Value = 0.3496932515337423
$ws.Range("B10").Value = 0.1123941493456505
$ws.Range("D10").Value = 0.02463433410315627
$ws.Range("E10").Value = 0.001539645881447267
$ws.Range("F10").Value = 0.06004618937644342
$ws.Range("J10").Value = 0.1293302540415704
$ws.Range("O10").Value = 0.01616628175519631
$ws.Range("Q10").Value = 0.262509622786759
$ws.Range("R10").Value = 0.05619707467282525
$ws.Range("S10").Value = 0.3371824480369515
$ws.Range("G11").Value = 0.1169811320754717
$ws.Range("J11").Value = 0.1169811320754717
$ws.Range("K11").Value = 0.169811320754717
$ws.Range("L11").Value = 0.5849056603773585
$ws.Range("S11").Value = 0.01132075471698113
$ws.Range("G12").Value = 0.6994535519125683
$ws.Range("J12").Value = 0.2076502732240437
$ws.Range("K12").Value = 0.01639344262295082
$ws.Range("L12").Value = 0.03278688524590164
$ws.Range("S12").Value = 0.04371584699453552
$ws.Range("F15").Value = 0.03482587064676617
$ws.Range("H15").Value = 0.1194029850746269
$ws.Range("I15").Value = 0.05970149253731343
$ws.Range("J15").Value = 0.3333333333333333
$ws.Range("K15").Value = 0.09950248756218906
$ws.Range("M15").Value = 0.004975124378109453
$ws.Range("O15").Value = 0.0945273631840796
$ws.Range("S15").Value = 0.2537313432835821
$ws.Range("F16").Value = 0.02051282051282051
$ws.Range("H16").Value = 0.2051282051282051
$ws.Range("I16").Value = 0.09743589743589744
$ws.Range("J16").Value = 0.3641025641025641
$ws.Range("K16").Value = 0.1076923076923077
$ws.Range("M16").Value = 0.02564102564102564
$ws.Range("S16").Value = 0.1128205128205128
$ws.Range("F17").Value = 0.01550387596899225
$ws.Range("H17").Value = 0.1744186046511628
$ws.Range("I17").Value = 0.05813953488372093
$ws.Range("J17").Value = 0.4476744186046512
$ws.Range("K17").Value = 0.08139534883720931
$ws.Range("M17").Value = 0.02713178294573643
$ws.Range("O17").Value = 0.07364341085271318
$ws.Range("S17").Value = 0.1220930232558139
$ws.Range("F18").Value = 0.0124223602484472
$ws.Range("H18").Value = 0.1739130434782609
$ws.Range("I18").Value = 0.08074534161490683
$ws.Range("J18").Value = 0.4347826086956522
$ws.Range("K18").Value = 0.09937888198757763
$ws.Range("M18").Value = 0.0124223602484472
$ws.Range("O18").Value = 0.03726708074534162
$ws.Range("S18").Value = 0.1490683229813665
$ws.Range("F19").Value = 0.0216076058772688
$ws.Range("H19").Value = 0.2221261884183232
$ws.Range("I19").Value = 0.07346585998271392
$ws.Range("J19").Value = 0.3837510803802939
$ws.Range("K19").Value = 0.09507346585998272
$ws.Range("M19").Value = 0.02074330164217805
$ws.Range("N19").Value = 0.000864304235090752
$ws.Range("O19").Value = 0.05704407951598963
$ws.Range("S19").Value = 0.125324114088159
